$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 164.63637
$ws.Range("I53").Value = 122.85714
$ws.Range("J53").Value = 237.75
$ws.Range("K53").Value = 122.85714
$ws.Range("L53").Value = 237.75
$ws.Range("M53").Value = 514.14286
$ws.Range("N53").Value = -1511.75

$ws.Range("H64").Value = 3412.75
$ws.Range("I64").Value = 3099.5
$ws.Range("J64").Value = 3517.1667
$ws.Range("K64").Value = 3099.5
$ws.Range("L64").Value = 3517.1667
$ws.Range("M64").Value = -2851.5
$ws.Range("N64").Value = -4013.1667

$ws.Range("H67").Value = 3412.75
$ws.Range("I67").Value = 3099.5
$ws.Range("J67").Value = 3517.1667
$ws.Range("K67").Value = 3099.5
$ws.Range("L67").Value = 3517.1667
$ws.Range("M67").Value = -2241.5
$ws.Range("N67").Value = -5233.1667

$ws.Range("H116").Value = 6492.4546
$ws.Range("I116").Value = 5717.8
$ws.Range("J116").Value = 7138
$ws.Range("K116").Value = 5717.8
$ws.Range("L116").Value = 7138
$ws.Range("M116").Value = -2275.8
$ws.Range("N116").Value = -14022

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6404.6816
$ws.Range("I63").Value = 5487.125
$ws.Range("K63").Value = 5487.125
$ws.Range("M63").Value = -4801.125

$ws.Range("H66").Value = 6404.6816
$ws.Range("I66").Value = 5487.125
$ws.Range("K66").Value = 27435.625
$ws.Range("M66").Value = -24003.625

$ws.Range("H122").Value = 387174.62
$ws.Range("I122").Value = 419022.5
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 1257067.5
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -1254617.5
$ws.Range("N122").Value = -19900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 651.3333
$ws.Range("I20").Value = 651.3333
$ws.Range("K20").Value = 651.3333
$ws.Range("M20").Value = -404.3333

$ws.Range("H96").Value = 11266.667
$ws.Range("I96").Value = 11266.667
$ws.Range("K96").Value = 11266.667
$ws.Range("M96").Value = -8520.666999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 83.40000000000001
$ws.Range("I7").Value = 66.14286
$ws.Range("J7").Value = 123.666664
$ws.Range("K7").Value = 66.14286
$ws.Range("L7").Value = 123.666664
$ws.Range("M7").Value = 46.85714
$ws.Range("N7").Value = -349.666664

$ws.Range("H86").Value = 9821
$ws.Range("I86").Value = 8878
$ws.Range("J86").Value = 10999.75
$ws.Range("K86").Value = 8878
$ws.Range("L86").Value = 10999.75
$ws.Range("M86").Value = -7755
$ws.Range("N86").Value = -13245.75

$ws.Range("H89").Value = 9821
$ws.Range("I89").Value = 8878
$ws.Range("J89").Value = 10999.75
$ws.Range("K89").Value = 44390
$ws.Range("L89").Value = 54998.75
$ws.Range("M89").Value = -38774
$ws.Range("N89").Value = -66230.75

$ws.Range("H122").Value = 4742.9443
$ws.Range("I122").Value = 4398.375
$ws.Range("J122").Value = 7499.5
$ws.Range("K122").Value = 13195.125
$ws.Range("L122").Value = 22498.5
$ws.Range("M122").Value = -10745.125
$ws.Range("N122").Value = -27398.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()

$ws.Range("H39").Value = 1107.5
$ws.Range("I39").Value = 730
$ws.Range("J39").Value = 2995
$ws.Range("K39").Value = 2190
$ws.Range("L39").Value = 8985
$ws.Range("M39").Value = -1896
$ws.Range("N39").Value = -9573

$ws.Range("H55").Value = 126848
$ws.Range("J55").Value = 4165
$ws.Range("L55").Value = 12495
$ws.Range("N55").Value = -12849

$ws.Range("H59").Value = 1500
$ws.Range("I59").Value = 1500
$ws.Range("K59").Value = 4500
$ws.Range("M59").Value = -3960

$ws.Range("H107").Value = 586.3143
$ws.Range("I107").Value = 225.66667
$ws.Range("J107").Value = 620.125
$ws.Range("K107").Value = 677.00001
$ws.Range("L107").Value = 1860.375
$ws.Range("M107").Value = 1242.99999
$ws.Range("N107").Value = -5700.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5499.5
$ws.Range("J70").Value = 4999.75
$ws.Range("L70").Value = 4999.75
$ws.Range("N70").Value = -5539.75

$ws.Range("H73").Value = 5499.5
$ws.Range("J73").Value = 4999.75
$ws.Range("L73").Value = 4999.75
$ws.Range("N73").Value = -6871.75

$ws.Range("H80").Value = 4146.1304
$ws.Range("I80").Value = 3394.4546
$ws.Range("J80").Value = 4835.1665
$ws.Range("K80").Value = 3394.4546
$ws.Range("L80").Value = 4835.1665
$ws.Range("M80").Value = -2396.4546
$ws.Range("N80").Value = -6831.1665

$ws.Range("H83").Value = 4146.1304
$ws.Range("I83").Value = 3394.4546
$ws.Range("J83").Value = 4835.1665
$ws.Range("K83").Value = 16972.273
$ws.Range("L83").Value = 24175.8325
$ws.Range("M83").Value = -11980.273
$ws.Range("N83").Value = -34159.8325

$ws.Range("H97").Value = 880.5
$ws.Range("I97").Value = 826.7
$ws.Range("K97").Value = 826.7
$ws.Range("M97").Value = -330.7

$ws.Range("H126").Value = 4267.933
$ws.Range("I126").Value = 3417.6667
$ws.Range("K126").Value = 10253.0001
$ws.Range("M126").Value = -7783.000100000001

$ws.Range("H132").Value = 4884.8
$ws.Range("I132").Value = 1198
$ws.Range("J132").Value = 7342.6665
$ws.Range("K132").Value = 3594
$ws.Range("L132").Value = 22027.9995
$ws.Range("M132").Value = -1064
$ws.Range("N132").Value = -27087.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4364
$ws.Range("I16").Value = 2182.2856
$ws.Range("J16").Value = 12000
$ws.Range("K16").Value = 2182.2856
$ws.Range("L16").Value = 12000
$ws.Range("M16").Value = -2012.2856
$ws.Range("N16").Value = -12340

$ws.Range("H40").Value = 1736.75
$ws.Range("I40").Value = 1736.75
$ws.Range("K40").Value = 1736.75
$ws.Range("M40").Value = -1600.75

$ws.Range("H46").Value = 1087.25
$ws.Range("I46").Value = 1383.3334
$ws.Range("K46").Value = 1383.3334
$ws.Range("M46").Value = -1195.3334

$ws.Range("H122").Value = 8055.5713
$ws.Range("I122").Value = 7077.8
$ws.Range("K122").Value = 21233.4
$ws.Range("M122").Value = -18783.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 572.1429000000001
$ws.Range("I113").Value = 715.8
$ws.Range("J113").Value = 213
$ws.Range("K113").Value = 2147.4
$ws.Range("L113").Value = 639
$ws.Range("M113").Value = 22.60000000000036
$ws.Range("N113").Value = -4979

$ws.Range("H122").Value = 6673052.5
$ws.Range("I122").Value = 6673052.5
$ws.Range("K122").Value = 20019157.5
$ws.Range("M122").Value = -20016707.5
